$d = $word.ActiveDocument

# --- Step 1: remove <w:rtl/> from the paragraph-mark properties of the
#     last paragraph (the "Disconnected domain model" answer paragraph)
#     by replacing its XML in place with the fixed version.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$fixedParaXml = @'
<w:p w14:paraId="6B00FE65" w14:textId="416381F2" w:rsidR="00CA02C2" w:rsidRPr="00EE7501" w:rsidRDefault="00CA02C2" w:rsidP="00CA02C2"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:bidi/><w:rPr><w:rFonts w:cs="B Nazanin"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">وقتی که داخل یه </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>aggregate</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> نیازه که با یه انتیتی یا </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>aggregate</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> دیگه ای کار کنیم و از </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>repository</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> اون انتیتی یا </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>aggregate</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> برای ارتباطش با </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>aggregate</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> مورد نظر استفاده می کنیم، </w:t></w:r><w:r w:rsidR="00307AFC"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">بهش گفته میشه </w:t></w:r><w:r w:rsidR="00307AFC"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>disconnected domain model</w:t></w:r><w:r w:rsidR="00307AFC"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>. چرا که دیگه خود اون</w:t></w:r><w:r w:rsidR="00307AFC" w:rsidRPr="00307AFC"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00307AFC"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">انتیتی یا </w:t></w:r><w:r w:rsidR="00307AFC"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>aggregate</w:t></w:r><w:r w:rsidR="00307AFC"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> داخل </w:t></w:r><w:r w:rsidR="00307AFC"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>domain</w:t></w:r><w:r w:rsidR="00307AFC"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> استفاده نشده.</w:t></w:r><w:r w:rsidR="00C4631F"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> چنین اتفاقی یه زنگ خطره به این معنی که احتمالا طراحی ما مشکل و نیاز به بازبینی داره.</w:t></w:r></w:p>
'@
$lastPara.Range.InsertXML($fixedParaXml)

# --- Step 2: append the new "second and third sections" content
#     (page break + the Eventual consistency / Domain event Q&A items)
#     at the very end of the document.
$endRange = $d.Content
$endRange.Collapse(0)
$newParasXml = @'
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:br w:type="page"/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:bidi/>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Eventual consistency</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> چیه؟</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:bidi/>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">در حالت عادی، توی یه برنامه وقتی یه درخواست سمت ما میاد، کل مراحلش توی یه بخش انجام میشه. در </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>eventual consistency</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> این مراحل پخش میشن و هر کدوم دست </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>domain</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> مخصوص خودشون میرن که انجام بشه و این پخش شدن از طریق </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>event</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> هاس. یه مثال بزنیم: تو همون بحث میتینگ که بالاتر گفتیم، اگه شخصی بخواد به میتینگ اضافه بشه، باید </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>domain</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> مربوط به میتینگ (مثلا </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>Meet</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">) چک کنه که آیا ظرفیت تکمیله یا نه، آیا این شخص اجازه اضافه شدن به میت رو داره یا نه و در کل مواردی که سمت </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>domain</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> مربوط به میتینگ باید چک بشه. از اونور باید بررسی بشه که مثلا </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ID</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> که از شخص گرفته شده درسته یا نه، شخص آنلاین هست یا نه و مواردی از این دست که مربوط به </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>domain</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> کاربر میشه. حالا اتفاقی که میفته اینه که هر بخش رو به </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>domain</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> خودش پاس میدیم که هندلش کنه و برای اینکار از </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>event</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> ها استفاده می کنیم.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:bidi/>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>Domain event</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> ها چی هستن؟</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:bidi/>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>یه تعریف باحال براش داریم:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">یه اتفاق </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>قابل اهمیت</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> (از دید بیزنس) که در سیستم به وقوع </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>پیوسته است</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:bidi/>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">خب اینکه یه تعریف خیلی کلی هست. اما در واقعیت این تعریف داره این ماجرا رو پوشش میده. داستان اینه که </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>یه سری کار ها هستن که برای ما مهمه وقتی که انجام میشن و به اتمام میرسن، بعدش یه کار دیگه ای انجام بشه. در واقع یه سری عملیات داریم که اتفاق افتادنشون به شرط اینه که یه عملیات دیگه ای قبلا به اتمام رسیده باشه.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> پس هر وقت که عملیات اصلی اتفاق افتاد و تموم شد، یه </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>event</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> رو </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>raise</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> می کنه که خبر بده آقا من کارم تموم شده حالا اونی که میخواسته بدونه، خبردار بشه تا کار خودشو انجام بده.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:bidi/>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">عملیاتی که میخواد بعد از اصلیه انجام بشه اما داخل </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>domain</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> مورد نظر هست معمولا. یعنی یه </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>domain</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> ای در واقع هست که اونه که میخواد خبردار بشه و یه کار دیگه ای رو انجام بده در نتیجه تعریف این </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>event</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> و پیاده سازیش میفته داخل </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>domain</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> مورد نظر و از این رو میشه یه </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>domain event</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin" w:hint="cs"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:bidi/>
        <w:rPr>
          <w:rFonts w:cs="B Nazanin"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@
$endRange.InsertXML($newParasXml)

Write-Output "Applied edits successfully"
